# Add newly uploaded daily COVID case rows (139-146, 2020-12-12 .. 2020-12-19)
# to the "cases_regions" sheet, following the exact same pattern used by the
# existing rows (e.g. row 138), and re-point the stray P135 "comment" cell
# at the normal comment style instead of the stale/duplicate one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B..M values (Total new, Adjara, Tbilisi, Imereti, Samegrelo-Zemo
# Svaneti, Kvemo Kartli, Shida Kartli, Kakheti, Guria, Mtskheta-Mtianeti,
# Samtskhe-Javakheti, Racha-Lechkhumi and Kvemo Svaneti) for every new row.
$rowData = @{
    139 = @(3907,333,1588,599,299,263,263,243,114,83,94,28)
    140 = @(2720,219,1240,362,203,231,156,142,40,84,43)
    141 = @(1337,212,571,147,74,168,48,41,16,41,18,1)
    142 = @(3837,263,1461,646,404,245,262,254,100,90,86,26)
    143 = @(3487,280,1477,533,269,254,283,153,66,85,76,11)
    144 = @(2981,228,1282,489,227,229,160,145,64,87,55,15)
    145 = @(2635,230,1134,338,199,216,169,186,66,46,39,12)
    146 = @(2904,197,1252,495,217,236,158,152,51,80,43,23)
}

# Column A date serials (1899 date system), one per new row.
$rowDates = @{
    139 = 44177
    140 = 44178
    141 = 44179
    142 = 44180
    143 = 44181
    144 = 44182
    145 = 44183
    146 = 44184
}

$firstRow = 139
$lastRow = 146

for ($r = $firstRow; $r -le $lastRow; $r++) {

    # Column A: date, formatted the same way as all the other date cells.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "YYYY\-MM\-DD"
    $dateCell.Value = $rowDates[$r]

    # Columns B..M: plain numeric values.
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }

    # Column O: the "Total new minus regional breakdown" check formula,
    # styled the same way as the rest of that column (e.g. O138).
    $oCell = $ws.Cells.Item($r, 15)
    $oCell.Style = $ws.Cells.Item($firstRow - 1, 15).Style
    $oCell.Formula = "=B" + $r + "-SUM(C" + $r + ":N" + $r + ")"
}

# P135 used to point at a stray/duplicate font-only style; repoint it to the
# same style used by the other "Comments" column cells (e.g. P124).
$ws.Cells.Item(135, 16).Style = $ws.Cells.Item(124, 16).Style

# Reflect where the user ended up after entering the new data.
$ws.Cells.Item(144, 15).Select()
